# "Bullet and Temp Level Objects added" - Will be textured shortly
#
# This script updates the "Todo list.xlsx" workbook:
#  - Programming sheet: updates progress status for "Floor Health" (Bullets
#    related work) and "Round Tracking" rows, and widens the Progress column
#    to fit the new, longer status text.
#  - Art sheet: marks the "Bullet Model" art task as Done (tackled by Zach),
#    and adds a new "Floor Layout" row describing a rough level-plan example.
#  - Selection/active sheet is moved to the Art sheet, matching where the
#    author ended up after making the edits.

$wb = $excel.ActiveWorkbook
$wsProgramming = $wb.Worksheets.Item("Programming")
$wsArt = $wb.Worksheets.Item("Art")

# ---------------------------------------------------------------------
# Programming sheet updates
# ---------------------------------------------------------------------

# "Floor Health" (row 8) - Bullets related floor-break code now needs
# artist support rather than being merely "In Progress".
$wsProgramming.Range("C8").Value = "Needs Artist Support"

# "Round Tracking" (row 13) - now done and can be linked into the UI.
$wsProgramming.Range("C13").Value = "Currently done and can be linked to UI"

# Widen the Progress column (C) so the longer status strings fit.
$wsProgramming.Columns.Item(3).ColumnWidth = 31.6

# Leave the selection on B6 (matches the author's final cursor position).
[void]$wsProgramming.Range("B6").Select()

# ---------------------------------------------------------------------
# Art sheet updates
# ---------------------------------------------------------------------

# "Bullet Model" (row 15) - Zach finished this, it's Done.
$wsArt.Range("B15").Value = "Zach"
$wsArt.Range("C15").Value = "Done"

# New row 18 - Floor Layout task, a rough example of the level plan.
$wsArt.Range("A18").Value = "Floor Layout"
$wsArt.Range("B18").Value = "Zach "
$wsArt.Range("C18").Value = "In Progress"
$wsArt.Range("E18").Value = "Rough example of the level plan"

# Make the Art sheet the active tab/sheet, with the selection on E18
# (matches the author's final cursor position).
[void]$wsArt.Activate()
[void]$wsArt.Range("E18").Select()
